# Automated update script applying the 2025-10-01 data refresh.
# Sheet 'VENTAS POR GRUPO': zero out a set of now-void group-sales cells
# and reset the 'N de 55' progress counters in row 57 to '0 de 55'.
# Sheet 'VENTA MENSUAL': roll the monthly columns forward by one month
# (junio/julio/agosto/septiembre -> julio/agosto/septiembre/octubre),
# shifting each client's monthly figures left accordingly, and adjust
# columns E/F widths to match the new month layout.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 'VENTAS POR GRUPO' -- zero out specific cells ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("C4").Value = 0
$ws1.Range("E4").Value = 0
$ws1.Range("L4").Value = 0
$ws1.Range("M4").Value = 0
$ws1.Range("D6").Value = 0
$ws1.Range("M6").Value = 0
$ws1.Range("D22").Value = 0
$ws1.Range("M22").Value = 0
$ws1.Range("D24").Value = 0
$ws1.Range("I25").Value = 0
$ws1.Range("M25").Value = 0
$ws1.Range("D26").Value = 0
$ws1.Range("D28").Value = 0
$ws1.Range("M28").Value = 0
$ws1.Range("E29").Value = 0
$ws1.Range("K29").Value = 0
$ws1.Range("M29").Value = 0
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 0
$ws1.Range("K31").Value = 0
$ws1.Range("L31").Value = 0
$ws1.Range("M31").Value = 0
$ws1.Range("C38").Value = 0
$ws1.Range("D38").Value = 0
$ws1.Range("E38").Value = 0
$ws1.Range("L38").Value = 0
$ws1.Range("M38").Value = 0
$ws1.Range("M39").Value = 0
$ws1.Range("K41").Value = 0
$ws1.Range("D44").Value = 0
$ws1.Range("H44").Value = 0
$ws1.Range("M44").Value = 0
$ws1.Range("D45").Value = 0
$ws1.Range("H45").Value = 0
$ws1.Range("I45").Value = 0
$ws1.Range("L45").Value = 0
$ws1.Range("M45").Value = 0
$ws1.Range("N45").Value = 0
$ws1.Range("L46").Value = 0
$ws1.Range("M46").Value = 0
$ws1.Range("I47").Value = 0
$ws1.Range("D49").Value = 0
$ws1.Range("M49").Value = 0
$ws1.Range("M52").Value = 0
$ws1.Range("N52").Value = 0
$ws1.Range("O53").Value = 0
$ws1.Range("L55").Value = 0
$ws1.Range("M55").Value = 0

# --- Sheet 1: row 57 summary labels, change leading count to 0 ---
$ws1.Range("C57").Value = "0 de 55"
$ws1.Range("D57").Value = "0 de 55"
$ws1.Range("E57").Value = "0 de 55"
$ws1.Range("H57").Value = "0 de 55"
$ws1.Range("I57").Value = "0 de 55"
$ws1.Range("K57").Value = "0 de 55"
$ws1.Range("L57").Value = "0 de 55"
$ws1.Range("M57").Value = "0 de 55"
$ws1.Range("N57").Value = "0 de 55"
$ws1.Range("O57").Value = "0 de 55"

# --- Sheet 2: 'VENTA MENSUAL' -- month header shift ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("C1").Value = "julio"
$ws2.Range("D1").Value = "agosto"
$ws2.Range("E1").Value = "septiembre"
$ws2.Range("F1").Value = "octubre"

# --- Sheet 2: shift monthly values left (C<-D, D<-E, E<-F, F<-0), with noted exceptions ---
$ws2.Range("E4").Value = 794.99
$ws2.Range("F4").Value = 0
$ws2.Range("C5").Value = 9991.16
$ws2.Range("D5").Value = 5652.48
$ws2.Range("E5").Value = 0
$ws2.Range("C6").Value = 1795.71
$ws2.Range("D6").Value = 1944.78
$ws2.Range("E6").Value = 1603.38
$ws2.Range("F6").Value = 0
$ws2.Range("C11").Value = 890.8
$ws2.Range("D11").Value = -295.8
$ws2.Range("E11").Value = 0
$ws2.Range("C12").Value = 320.98
$ws2.Range("D12").Value = 0
$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 174.18
$ws2.Range("E14").Value = 0
$ws2.Range("C20").Value = 0
$ws2.Range("C22").Value = 1710.72
$ws2.Range("D22").Value = 8385.889999999999
$ws2.Range("E22").Value = 146.99
$ws2.Range("F22").Value = 0
$ws2.Range("E24").Value = 457.92
$ws2.Range("F24").Value = 0
$ws2.Range("C25").Value = 818.0599999999999
$ws2.Range("D25").Value = 497.66
$ws2.Range("E25").Value = 8255.23
$ws2.Range("F25").Value = 0
$ws2.Range("C26").Value = 1373.76
$ws2.Range("D26").Value = 475.2
$ws2.Range("E26").Value = 950.4
$ws2.Range("F26").Value = 0
$ws2.Range("C27").Value = 0
$ws2.Range("C28").Value = 17469.82
$ws2.Range("D28").Value = 9158.4
$ws2.Range("E28").Value = 10350.26
$ws2.Range("F28").Value = 0
$ws2.Range("C29").Value = 14529.6
$ws2.Range("D29").Value = 222.32
$ws2.Range("E29").Value = 8733.540000000001
$ws2.Range("F29").Value = 0
$ws2.Range("D30").Value = 61.75
$ws2.Range("E30").Value = 0
$ws2.Range("C31").Value = 5996.2
$ws2.Range("D31").Value = 486.71
$ws2.Range("E31").Value = 10174.33
$ws2.Range("F31").Value = 0
$ws2.Range("C32").Value = 739.1
$ws2.Range("D32").Value = 0
$ws2.Range("C33").Value = 0
$ws2.Range("D33").Value = 518.4
$ws2.Range("E33").Value = 0
$ws2.Range("C38").Value = 4677
$ws2.Range("D38").Value = 7942.96
$ws2.Range("E38").Value = 8322.860000000001
$ws2.Range("F38").Value = 0
$ws2.Range("C39").Value = 1186.08
$ws2.Range("D39").Value = 0
$ws2.Range("E39").Value = 1428.84
$ws2.Range("F39").Value = 0
$ws2.Range("C41").Value = 2874.67
$ws2.Range("D41").Value = 660.24
$ws2.Range("E41").Value = 7011.36
$ws2.Range("F41").Value = 0
$ws2.Range("C43").Value = 582.66
$ws2.Range("D43").Value = 295.63
$ws2.Range("E43").Value = 0
$ws2.Range("C44").Value = 3269.52
$ws2.Range("D44").Value = 1824.2
$ws2.Range("E44").Value = 1949.23
$ws2.Range("F44").Value = 0
$ws2.Range("C45").Value = 5003.99
$ws2.Range("D45").Value = 1151.68
$ws2.Range("E45").Value = 4381.83
$ws2.Range("F45").Value = 0
$ws2.Range("C46").Value = 6207.1
$ws2.Range("D46").Value = 11132.19
$ws2.Range("E46").Value = 2781.14
$ws2.Range("F46").Value = 0
$ws2.Range("C47").Value = 1712.88
$ws2.Range("D47").Value = 1549.1
$ws2.Range("E47").Value = 71.62
$ws2.Range("F47").Value = 0
$ws2.Range("C49").Value = 154.28
$ws2.Range("D49").Value = 6923.38
$ws2.Range("E49").Value = 108.12
$ws2.Range("F49").Value = 0
$ws2.Range("C51").Value = 380.16
$ws2.Range("D51").Value = 0
$ws2.Range("C52").Value = 1439.92
$ws2.Range("D52").Value = 0
$ws2.Range("E52").Value = 3995.75
$ws2.Range("F52").Value = 0
$ws2.Range("E53").Value = 869.53
$ws2.Range("F53").Value = 0
$ws2.Range("D54").Value = 142.2
$ws2.Range("E54").Value = 0
$ws2.Range("C55").Value = 581.26
$ws2.Range("D55").Value = 0
$ws2.Range("E55").Value = 1556.27
$ws2.Range("F55").Value = 0
$ws2.Range("C57").Value = 83705.42999999999
$ws2.Range("D57").Value = 58903.55
$ws2.Range("E57").Value = 73943.59
$ws2.Range("F57").Value = 0

# --- Sheet 2: adjust column widths for E (->16) and F (->13) ---
$ws2.Columns.Item(5).ColumnWidth = 15.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 12.166666666666666
